$wb = $excel.ActiveWorkbook

# ALC!row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5371.143
$ws.Range("I74").Value = 3799
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 3799
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -2863
$ws.Range("N74").Value = -7872

# ALC!row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 5371.143
$ws.Range("I77").Value = 3799
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 18995
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -14315
$ws.Range("N77").Value = -39360

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1569.55
$ws.Range("I137").Value = 1749.25
$ws.Range("J137").Value = 1300
$ws.Range("K137").Value = 5247.75
$ws.Range("L137").Value = 3900
$ws.Range("M137").Value = -2697.75
$ws.Range("N137").Value = -9000

# ARM!row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25647460
$ws.Range("I32").Value = 5582.1113
$ws.Range("J32").Value = 333350000
$ws.Range("K32").Value = 5582.1113
$ws.Range("L32").Value = 333350000
$ws.Range("M32").Value = -5295.1113
$ws.Range("N32").Value = -333350574

# ARM!row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2370
$ws.Range("I45").Value = 2000
$ws.Range("J45").Value = 2555
$ws.Range("K45").Value = 2000
$ws.Range("L45").Value = 2555
$ws.Range("M45").Value = -1623
$ws.Range("N45").Value = -3309

# ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13890387
$ws.Range("I61").Value = 13890387
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 13890387
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -13890175
$ws.Range("N61").ClearContents()

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1279882.6
$ws.Range("I132").Value = 920.1622
$ws.Range("J132").Value = 6537839.5
$ws.Range("K132").Value = 2760.4866
$ws.Range("L132").Value = 19613518.5
$ws.Range("M132").Value = -230.4866000000002
$ws.Range("N132").Value = -19618578.5

# ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 13890387
$ws.Range("I136").Value = 13890387
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 41671161
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -41668611
$ws.Range("N136").ClearContents()

# BSM!row 5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1308.4
$ws.Range("I5").Value = 18.5
$ws.Range("J5").Value = 2168.3333
$ws.Range("K5").Value = 18.5
$ws.Range("L5").Value = 2168.3333
$ws.Range("M5").Value = 94.5
$ws.Range("N5").Value = -2394.3333

# BSM!row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 55557644
$ws.Range("I107").Value = 83334460
$ws.Range("J107").Value = 3993.3333
$ws.Range("K107").Value = 83334460
$ws.Range("L107").Value = 3993.3333
$ws.Range("M107").Value = -83332540
$ws.Range("N107").Value = -7833.3333

# BSM!row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7945930
$ws.Range("I134").Value = 2913.875
$ws.Range("J134").Value = 18536618
$ws.Range("K134").Value = 8741.625
$ws.Range("L134").Value = 55609854
$ws.Range("M134").Value = -6206.625
$ws.Range("N134").Value = -55614924

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1390097.2
$ws.Range("I31").Value = 1792954.5
$ws.Range("J31").Value = 2477.7778
$ws.Range("K31").Value = 1792954.5
$ws.Range("L31").Value = 2477.7778
$ws.Range("M31").Value = -1792659.5
$ws.Range("N31").Value = -3067.7778

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1390097.2
$ws.Range("I34").Value = 1792954.5
$ws.Range("J34").Value = 2477.7778
$ws.Range("K34").Value = 1792954.5
$ws.Range("L34").Value = 2477.7778
$ws.Range("M34").Value = -1792752.5
$ws.Range("N34").Value = -2881.7778

# CUL!row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 36462080
$ws.Range("I5").Value = 50724980
$ws.Range("J5").Value = 12444.444
$ws.Range("K5").Value = 152174940
$ws.Range("L5").Value = 37333.33199999999
$ws.Range("M5").Value = -152174828
$ws.Range("N5").Value = -37557.33199999999

# CUL!row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 36462080
$ws.Range("I135").Value = 50724980
$ws.Range("J135").Value = 12444.444
$ws.Range("K135").Value = 456524820
$ws.Range("L135").Value = 111999.996
$ws.Range("M135").Value = -456522285
$ws.Range("N135").Value = -117069.996

# CUL!row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 10418418
$ws.Range("I140").Value = 13890401
$ws.Range("K140").Value = 41671203
$ws.Range("M140").Value = -41666023

# GSM!row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5265560
$ws.Range("I80").Value = 2373.182
$ws.Range("J80").Value = 12502442
$ws.Range("K80").Value = 2373.182
$ws.Range("L80").Value = 12502442
$ws.Range("M80").Value = -1375.182
$ws.Range("N80").Value = -12504438

# GSM!row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 5265560
$ws.Range("I83").Value = 2373.182
$ws.Range("J83").Value = 12502442
$ws.Range("K83").Value = 11865.91
$ws.Range("L83").Value = 62512210
$ws.Range("M83").Value = -6873.91
$ws.Range("N83").Value = -62522194

# GSM!row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2148.5
$ws.Range("I126").Value = 1596
$ws.Range("J126").Value = 2332.6667
$ws.Range("K126").Value = 4788
$ws.Range("L126").Value = 6998.000100000001
$ws.Range("M126").Value = -2318
$ws.Range("N126").Value = -11938.0001

# GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5583.963
$ws.Range("I132").Value = 1480.1177
$ws.Range("J132").Value = 12560.5
$ws.Range("K132").Value = 4440.3531
$ws.Range("L132").Value = 37681.5
$ws.Range("M132").Value = -1910.3531
$ws.Range("N132").Value = -42741.5

# LTW!row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 45456844
$ws.Range("I61").Value = 2227
$ws.Range("J61").Value = 166669150
$ws.Range("K61").Value = 2227
$ws.Range("L61").Value = 166669150
$ws.Range("M61").Value = -2025
$ws.Range("N61").Value = -166669554

# LTW!row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 45456844
$ws.Range("I113").Value = 2227
$ws.Range("J113").Value = 166669150
$ws.Range("K113").Value = 2227
$ws.Range("L113").Value = 166669150
$ws.Range("M113").Value = -57
$ws.Range("N113").Value = -166673490

# LTW!row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 43966104
$ws.Range("I132").Value = 87913130
$ws.Range("J132").Value = 19075.924
$ws.Range("K132").Value = 263739390
$ws.Range("L132").Value = 57227.772
$ws.Range("M132").Value = -263736860
$ws.Range("N132").Value = -62287.772

# WVR!row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1109.7778
$ws.Range("I107").Value = 798
$ws.Range("J107").Value = 1359.2
$ws.Range("K107").Value = 2394
$ws.Range("L107").Value = 4077.6
$ws.Range("M107").Value = -474
$ws.Range("N107").Value = -7917.6

# WVR!row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1269.4615
$ws.Range("I136").Value = 786.17645
$ws.Range("K136").Value = 2358.52935
$ws.Range("M136").Value = 191.4706499999998
